$d = $word.ActiveDocument

# --- 1. Merge the two identically-formatted runs in the "In our project..." paragraph
#        into a single run by touching the combined range (Word recombines adjacent
#        runs that share formatting whenever the text in that span is rewritten).
$mergedText = "In our project, we implement interceptor design pattern in logger operation and we use priority strategy as callback strategy. As our project designed, each interceptor has its own priority, and the dispatcher will always invoke the interceptor who has the most large priority number currently in the list."
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)

# --- 2. Drop the old "_GoBack" bookmark that sits after "Add value priority of interceptor".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Insert a new heading paragraph at the very top of the document.
$firstPara = $d.Paragraphs(1).Range
$firstPara.InsertParagraphBefore()
$newPara = $d.Paragraphs(1).Range
$newPara.Text = "New version with more code"

# --- 4. Re-create "_GoBack" so it wraps the freshly typed text (bookmarkStart right
#        before the run, bookmarkEnd right after it, matching the original placement
#        style used elsewhere in this document).
$bmRange = $d.Range($newPara.Start, $newPara.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
